# Edit slide 1 / "TextBox 19" (shape 8):
#   1. "Hackidemy" -> "Hackademy"
#   2. " is a good learning program. " split into " " + "is a good learning program. "
#      (same wording, just now two separate runs instead of one)
#   3. " to change the " + "world" merged back into a single run " to change the world"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(8)
$tr  = $shp.TextFrame.TextRange

# 1) Fix the typo: Hackidemy -> Hackademy
$hit = $tr.Find("Hackidemy")
$hit.Text = "Hackademy"

# 2) Split " is a good learning program. " into two runs: " " and "is a good learning program. "
$hit = $tr.Find(" is a good learning program. ")
$start = $hit.Start
$tr.Characters($start, 1).Text = " "
$tr.Characters($start + 1, 28).Text = "is a good learning program. "

# 3) Merge " to change the " + "world" into a single run " to change the world"
$hit = $tr.Find(" to change the ")
$start = $hit.Start
$tr.Characters($start, 20).Text = " to change the world"
